$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers (row 1) and values (row 2) for the two additional columns
$ws.Range("G1").Value = "Creation monstre et equipement"
$ws.Range("H1").Value = "Gestion equipement"

$ws.Range("G2").Value = "Yoan"
$ws.Range("H2").Value = "Enzo"

# Match the style used by the existing header row (wrap text)
$ws.Range("G1:H1").WrapText = $true

# Update selection to mirror the author's final cursor position
$ws.Range("H3").Select()
